$wb = $excel.ActiveWorkbook

# --- Reorder sheets: "produtos" should come before "busca_cesta" ---
$wb.Worksheets("produtos").Move($wb.Worksheets("busca_cesta"))

$wsP = $wb.Worksheets("produtos")
$wsB = $wb.Worksheets("busca_cesta")

# --- Remove the "na_cesta" (column C) marker for these 10 products ---
# (they were removed from the current cesta/basket)
$rowsToClear = @(5, 12, 15, 19, 24, 31, 39, 43, 48, 53)
foreach ($r in $rowsToClear) {
    $wsP.Range("C$r").ClearContents()
}

# --- Update view/selection state on busca_cesta ---
$wsB.Activate()
$wsB.Range("I19").Select()

# --- Update view/selection state on produtos (now the active sheet) ---
$wsP.Activate()
$wsP.Range("C2:C55").Select()
